# CR Upload documents -4/4
# - Fix typo "Uplode" -> "Upload" in the CR sheet's General notes header/column
#   (propagated down the whole F1:F10 column, matching the other column
#   values that repeat per-row).
# - Add a new "Other Supporting documents" column (G1:G10).
# - Update the selection / used-range to include the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CR")

# Column F: corrected "Change request Upload document General notes" text,
# repeated for every data row (1-10).
$ws.Range("F1:F10").Value = "Change request Upload document General notes"

# Column G: new "Other Supporting documents" column, same value for every row.
$ws.Range("G1:G10").Value = "Other Supporting documents"

# Best-fit-ish widths for the new/changed columns so the sheet reads cleanly.
$ws.Columns.Item(6).ColumnWidth = 44.25

# Match the saved selection: G1:G10 with the active cell at the top (G1).
[void]$ws.Range("G1:G10").Select()
